# Auto-generated script applying 2025-03-15 daily crime data updates
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 1115
$ws.Range("L3").Value = 1127
$ws.Range("F4").Value = 1923
$ws.Range("J4").Value = 1857
$ws.Range("K4").Value = 1746
$ws.Range("L4").Value = 314
$ws.Range("K6").Value = 9122
$ws.Range("L6").Value = 1146
$ws.Range("F7").Value = 24116
$ws.Range("J7").Value = 29329
$ws.Range("K7").Value = 27538
$ws.Range("L7").Value = 3777

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 69
$ws.Range("L6").Value = 72
$ws.Range("L7").Value = 235

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 25
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L3").Value = 33
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 24
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 117
$ws.Range("L8").Value = 235
$ws.Range("L9").Value = 22
$ws.Range("L11").Value = 61
$ws.Range("L14").Value = 16
$ws.Range("L15").Value = 26
$ws.Range("L19").Value = 117
$ws.Range("L20").Value = 100
$ws.Range("K27").Value = 260
$ws.Range("L27").Value = 43
$ws.Range("L29").Value = 176
$ws.Range("L35").Value = 5
$ws.Range("L37").Value = 127
$ws.Range("L42").Value = 118
$ws.Range("L43").Value = 29
$ws.Range("L47").Value = 32
$ws.Range("L51").Value = 50
$ws.Range("L52").Value = 70
$ws.Range("L53").Value = 46
$ws.Range("L60").Value = 22
$ws.Range("F63").Value = 208
$ws.Range("J63").Value = 206
$ws.Range("K63").Value = 83
$ws.Range("L63").Value = 17
$ws.Range("L65").Value = 75
$ws.Range("L66").Value = 5
$ws.Range("L67").Value = 140
$ws.Range("L69").Value = 9
$ws.Range("L70").Value = 15
$ws.Range("L72").Value = 15
$ws.Range("L73").Value = 28
$ws.Range("L76").Value = 45
$ws.Range("L79").Value = 106
$ws.Range("L83").Value = 88
$ws.Range("L85").Value = 195
$ws.Range("L87").Value = 14
$ws.Range("L91").Value = 51
$ws.Range("L95").Value = 56
$ws.Range("L97").Value = 46
$ws.Range("L98").Value = 31
$ws.Range("L99").Value = 59
$ws.Range("F101").Value = 24116
$ws.Range("J101").Value = 29329
$ws.Range("K101").Value = 27538
$ws.Range("L101").Value = 3777

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 46
$ws.Range("L3").Value = 38
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 62
$ws.Range("L3").Value = 58
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 176

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 33
$ws.Range("L3").Value = 39
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 30
$ws.Range("L4").Value = 9
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 13
$ws.Range("L6").Value = 20

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 34
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L2").Value = 14
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("L3").Value = 2

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("L7").Value = 5

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("L3").Value = 1
$ws.Range("L7").Value = 5

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L2").Value = 10
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L3").Value = 15
$ws.Range("K4").Value = 32
$ws.Range("K7").Value = 260
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 85
$ws.Range("L4").Value = 16
$ws.Range("L6").Value = 41
$ws.Range("L7").Value = 195

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 24
$ws.Range("L3").Value = 19
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("L3").Value = 1
$ws.Range("L7").Value = 14
